# Applies the diff described:
#  - software sheet: update report_date (B1)
#  - qc sheet: update row 6 (C6, D6, E6)
#  - emu_long sheet: update rows 12-21 (barcode04 / barcode05 breakdown),
#    including swapping row 17/18 content + styling row 17 blue (s=3)
#  - emu_counts sheet: update L6, L9
#  - emu_proportions sheet: update K2, K4, K5, L6, L7, L8, K9, L9

$wb = $excel.ActiveWorkbook

# ---- software sheet ----
$wsSoftware = $wb.Worksheets.Item("software")
$wsSoftware.Range("B1").Value = 45758.65994086969

# ---- qc sheet ----
$wsQc = $wb.Worksheets.Item("qc")
$wsQc.Range("C6").Value = 55
$wsQc.Range("D6").Value = 50
$wsQc.Range("E6").Value = 0.5

# ---- emu_long sheet ----
$wsLong = $wb.Worksheets.Item("emu_long")

# row 12 (barcode04 / unassigned)
$wsLong.Range("C12").Value = 0

# row 13 (barcode04 / Escherichia coli)
$wsLong.Range("C13").Value = 0.74

# row 14 (barcode04 / Staphylococcus epidermidis)
$wsLong.Range("C14").Value = 0.19

# row 15 (barcode04 / Staphylococcus aureus)
$wsLong.Range("C15").Value = 0.07000000000000001

# row 17 becomes Leucobacter aridicollis (was row 18's species), styled blue (s=3)
$wsLong.Range("A17:F17").Style = "Normal"
$wsLong.Range("B17").Value = "Leucobacter aridicollis"
$wsLong.Range("C17").Value = 0.91
$wsLong.Range("D17").Value = 50
$wsLong.Range("E17").Value = 0.48
$wsLong.Range("F17").Value = 47.62
$wsLong.Range("A17:F17").Font.Color = 16711680

# row 18 becomes unassigned (was row 17's species)
$wsLong.Range("B18").Value = "unassigned"
$wsLong.Range("C18").Value = 0
$wsLong.Range("D18").Value = 50
$wsLong.Range("E18").Value = 0.48
$wsLong.Range("F18").Value = 47.62

# row 19 (barcode05 / Bacillus sp. IHB B 7164)
$wsLong.Range("C19").Value = 0.05
$wsLong.Range("F19").Value = 2.64

# row 20 (barcode05 / Bacillus megaterium)
$wsLong.Range("C20").Value = 0.04
$wsLong.Range("F20").Value = 2.12

# row 21 (barcode05 / total)
$wsLong.Range("D21").Value = 105

# ---- emu_counts sheet ----
$wsCounts = $wb.Worksheets.Item("emu_counts")
$wsCounts.Range("L6").Value = 50
$wsCounts.Range("L9").Value = 50

# ---- emu_proportions sheet ----
$wsProp = $wb.Worksheets.Item("emu_proportions")
$wsProp.Range("K2").Value = 0.19
$wsProp.Range("K4").Value = 0.07000000000000001
$wsProp.Range("K5").Value = 0.74
$wsProp.Range("L6").Value = 0.91
$wsProp.Range("L7").Value = 0.05
$wsProp.Range("L8").Value = 0.04
$wsProp.Range("K9").Value = 0
$wsProp.Range("L9").Value = 0
